$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Update the page count (bookmark) for "Researching Information Systems and Computing" (row 11)
$ws.Range("C11").Value = 219

# Widen column B slightly to fit the new, longer author/title text
$ws.Columns.Item(2).ColumnWidth = 31.5

# Add the new book entry in row 16
$ws.Range("A16").Value = "Artificial Intelligence: Foundations of Computational Agents (2nd Edition)"
$ws.Range("B16").Value = "David L. Poole and Alan K. Mackworth"

# Update the selection to reflect where the user was last working
$ws.Range("C16").Select()

$wb.Save()
